# Generate Report for Handback
# Update handback status timestamps / status for the a811ebfb... and
# e6c0ce55... rows (rows 3 and 5 across the Overview / zh-cn / de-de sheets),
# which previously shared identical string values and therefore move
# together.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G3").Value = "2016-08-12 14:15:58"
$wsOverview.Range("G5").Value = "2016-08-12 14:15:58"

# zh-cn sheet - Status (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K)
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-12 14:15:51"
$wsZhCn.Range("H5").Value = "2016-08-12 14:15:51"
$wsZhCn.Range("K3").Value = "2016-08-12 14:16:22"
$wsZhCn.Range("K5").Value = "2016-08-12 14:16:22"

# de-de sheet - Status (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K)
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-12 14:15:58"
$wsDeDe.Range("H5").Value = "2016-08-12 14:15:58"
$wsDeDe.Range("K3").Value = "2016-08-12 14:16:32"
$wsDeDe.Range("K5").Value = "2016-08-12 14:16:32"
